# Se corrige reporte de factura, remitos. Se reescriben consultas de CC de
# clientes y proveedores.
#
# Adds two new backlog rows to the "Hoja1" task list and leaves the
# selection/scroll position where the user ended up after typing them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Scroll the view down towards the bottom of the (growing) list, like a user
# paging down to add new rows below the existing ones.
$excel.Goto($ws.Range("A22"))

# New backlog item: facturación de remitos pendientes en dólares.
$ws.Range("A34").Value = "facturacion de remitos pendientes en dolares esta fallando"
$ws.Range("B34").Value = "terminado"

# New backlog item: balance en cuenta corriente de proveedores.
$ws.Range("A35").Value = "balance en cc proveedores esta fallando"
$ws.Range("B35").Value = "terminado"

# Leave the selection on the next empty row, as Excel would after typing the
# last entry and pressing Enter/Tab into the following blank row.
$ws.Range("B37").Select()
